$p = $ppt.ActivePresentation

# --- locate the slide / shape / run that holds the sentence we need to edit ---
$oldSentence = "Eliminates occurrence of a word(cause) from positive sentiment-sorted array "
$newLead     = "Eliminates occurrence of a word(cause) in positive sentiment array from negative sentiment array and "
$newTrail    = "vice versa. "

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -like "*$oldSentence*") {
                $targetSlide = $sl
                $targetShape = $sh
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf($oldSentence)

# --- 1) rewrite the existing run's text (keeps its original rPr / dirty="0") ---
$run1 = $tr.Characters($idx + 1, $oldSentence.Length)
$run1.Text = $newLead

# --- 2) type the trailing sentence right after it as a brand-new run ---
$run1Fresh = $tr.Characters($idx + 1, $newLead.Length)
$run1Fresh.InsertAfter($newTrail) | Out-Null

# Force the newly-typed text to stay a distinct <a:r> (matching its own
# explicit bold/size/color formatting) instead of silently re-merging with
# the previous run.
$run2 = $tr.Characters($idx + 1 + $newLead.Length, $newTrail.Length)
$run2.Font.Bold = $true
$run2.Font.Size = 14
